# NN-368 - Paymaart - Admin Web - Insights of Merchant Registration BDD issue solved.
# Update the sample "Merchant Registration" row on Sheet1 with corrected
# test fixture values (Paymaart ID + Phone Number), matching the fix for
# the BDD test data, and nudge the header/data row heights back up to the
# sheet's default so row 1-3 render consistently with the rest of the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = "Paymaart ID", Column D = "Phone Number" (row 2 = first data row)
$ws.Range("C2").Value = "CMR47857280"
$ws.Range("D2").Value = "265 84 546 3484 "

# Restore header/data rows to the sheet's normal row height.
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# Leave the cursor on the phone number cell that was just edited.
$ws.Range("D2").Select()
